$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 3.25
$ws.Range("J2").Value = 4.33
$ws.Range("K2").Value = 1.95
$ws.Range("L2").Value = 3
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
$ws.Range("S2").Value = 1.57
$ws.Range("T2").Value = 2.25
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 8
$ws.Range("AB2").Value = 41
$ws.Range("AC2").Value = 7
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 6
$ws.Range("AI2").Value = 9.5
$ws.Range("AJ2").Value = 19
$ws.Range("AM2").Value = 501
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 34
$ws.Range("AS2").Value = 351
$ws.Range("AT2").Value = 2.25
$ws.Range("AU2").Value = 9
$ws.Range("AV2").Value = 67
$ws.Range("AZ2").Value = 41
$ws.Range("BB2").Value = 251
$ws.Range("BD2").Value = 151

# Row 8
$ws.Range("K8").Value = 2.38
$ws.Range("W8").Value = 6.5
$ws.Range("Z8").Value = 8.5
$ws.Range("AW8").Value = 9

# Row 10
$ws.Range("G10").Value = 2.4
$ws.Range("I10").Value = 3.2
$ws.Range("Y10").Value = 10
$ws.Range("AH10").Value = 15
$ws.Range("AI10").Value = 13
$ws.Range("AK10").Value = 34

# Row 16
$ws.Range("AU16").Value = 7.5
$ws.Range("AY16").Value = 21
$ws.Range("BA16").Value = 51
